# Rename sheets per new naming scheme (rerun LU d2c FeatEng for FR cities)
$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ57551917",
    "summ57721554",
    "summ57904853",
    "summ58121322",
    "summ58347551",
    "summ58572540",
    "summ58804759",
    "summ59038589",
    "summ59271948",
    "summ59505301",
    "summ59740591",
    "summ59971333",
    "summ00249074",
    "summ00471355",
    "summ00704876",
    "summ00938036",
    "summ01176235",
    "summ01409679",
    "summ01638110",
    "summ01872052",
    "summ02105008",
    "summ02338055",
    "summ02554668",
    "summ02788864",
    "summ03023963",
    "summ03265733",
    "summ03488091",
    "summ03721845",
    "summ03953203",
    "summ04171270",
    "summ04405622",
    "summ04638031",
    "summ04881727",
    "summ05123309",
    "summ05360445",
    "summ05610711",
    "summ05856854",
    "summ06088796",
    "summ06328033",
    "summ06555380",
    "summ06800615",
    "summ07021920",
    "summ07271935",
    "summ07488724",
    "summ07727954",
    "summ07973587",
    "summ08222221",
    "summ08454697",
    "summ08671325",
    "summ08921264",
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

